# Updates cryptos list price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.238.02'
$ws.Range("E2").Value = '  +0.29%  '
# Row 3
$ws.Range("D3").Value = '2.305.88'
$ws.Range("E3").Value = '  -0.06%  '
# Row 4
$ws.Range("E4").Value = '  +0.09%  '
# Row 5
$ws.Range("D5").Value = '''302.26'
$ws.Range("E5").Value = '  +0.32%  '
# Row 6
$ws.Range("D6").Value = '''100.17'
$ws.Range("E6").Value = '  +2.95%  '
# Row 7
$ws.Range("E7").Value = '  +0.48%  '
# Row 8
$ws.Range("E8").Value = '  +0.08%  '
# Row 9
$ws.Range("E9").Value = '  +3.81%  '
# Row 10
$ws.Range("D10").Value = '''36.36'
$ws.Range("E10").Value = '  +8.11%  '
# Row 11
$ws.Range("E11").Value = '  +0.08%  '
# Row 12
$ws.Range("D12").Value = '''18.68'
$ws.Range("E12").Value = '  +10.31%  '
# Row 13
$ws.Range("E13").Value = '  +1.01%  '
# Row 14
$ws.Range("E14").Value = '  +3.01%  '
# Row 15
$ws.Range("D15").Value = '2.667.37'
$ws.Range("E15").Value = '  +0.35%  '
# Row 16
$ws.Range("D16").Value = '2.313.18'
$ws.Range("E16").Value = '  +0.13%  '
# Row 17
$ws.Range("E17").Value = '  -0.03%  '
# Row 18
$ws.Range("D18").Value = '43.110.59'
$ws.Range("E18").Value = '  +0.51%  '
# Row 19
$ws.Range("D19").Value = '''12.90'
$ws.Range("E19").Value = '  +10.98%  '
# Row 20
$ws.Range("D20").Value = '''6.18'
$ws.Range("E20").Value = '  +2.83%  '
# Row 21
$ws.Range("D21").Value = '0.0₃0907'
$ws.Range("E21").Value = '  +0.73%  '
# Row 22
$ws.Range("D22").Value = '''68.11'
$ws.Range("E22").Value = '  +1.45%  '
# Row 23
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").Value = '''2.26'
$ws.Range("E23").Value = '  +13.95%  '
# Row 24
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '''236.56'
$ws.Range("E24").Value = '  -0.06%  '
# Row 25
$ws.Range("D25").Value = '''1.01'
$ws.Range("E25").Value = '  +0.44%  '
# Row 26
$ws.Range("E26").Value = '  -0.05%  '
# Row 27
$ws.Range("E27").Value = '  +1.19%  '
# Row 28
$ws.Range("E28").Value = '  +8.94%  '
# Row 29
$ws.Range("D29").Value = '''34.91'
$ws.Range("E29").Value = '  +2.75%  '
# Row 30
$ws.Range("D30").Value = '''167.50'
$ws.Range("E30").Value = '  +0.68%  '
# Row 31
$ws.Range("D31").Value = '''9.17'
$ws.Range("E31").Value = '  +0.42%  '
# Row 32
$ws.Range("D32").Value = '''0.999'
$ws.Range("E32").Value = '  +0.02%  '
# Row 33
$ws.Range("D33").Value = '''5.05'
$ws.Range("E33").Value = '  +1.58%  '
# Row 34
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = '''4.73'
$ws.Range("E34").Value = '  -1.01%  '
# Row 35
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").Value = '''17.88'
$ws.Range("E35").Value = '  +5.30%  '
# Row 36
$ws.Range("E36").Value = '  +0.82%  '
# Row 37
$ws.Range("D37").Value = '''0.0697'
$ws.Range("E37").Value = '  -0.18%  '
# Row 38
$ws.Range("E38").Value = '  +0.37%  '
# Row 39
$ws.Range("E39").Value = '  +2.33%  '
# Row 40
$ws.Range("E40").Value = '  +0.67%  '
# Row 41
$ws.Range("E41").Value = '  +0.74%  '
# Row 42
$ws.Range("D42").Value = '''2.35'
$ws.Range("E42").Value = '  +1.06%  '
# Row 43
$ws.Range("D43").Value = '1.991.70'
$ws.Range("E43").Value = '  +1.13%  '
# Row 44
$ws.Range("E44").Value = '  +3.52%  '
# Row 45
$ws.Range("D45").Value = '''10.14'
$ws.Range("E45").Value = '  +3.66%  '
# Row 46
$ws.Range("D46").Value = '''17.72'
$ws.Range("E46").Value = '  +0.18%  '
# Row 47
$ws.Range("D47").Value = '''2.90'
$ws.Range("E47").Value = '  +1.73%  '
# Row 48
$ws.Range("D48").Value = '''55.72'
$ws.Range("E48").Value = '  +5.38%  '
# Row 49
$ws.Range("E49").Value = '  +3.99%  '
# Row 50
$ws.Range("D50").Value = '2.530.26'
$ws.Range("E50").Value = '  +0.26%  '
# Row 51
$ws.Range("D51").Value = '''70.96'
$ws.Range("E51").Value = '  +1.46%  '
